$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.918.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.339.05'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.20%  '
$ws.Range("E4").Value = '  -0.67%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.577'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.06%  '
$ws.Range("E11").Value = '  +1.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.43%  '
$ws.Range("E13").Value = '  -0.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.692.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.333.66'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.15'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.63%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.829'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '46.778.84'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +15.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0949'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '246.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.54%  '
$ws.Range("E24").Value = '  +1.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '42.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +15.02%  '
$ws.Range("E28").Value = '  +0.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.90'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.78'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.37'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0819'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.07%  '
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.18'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.113'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.01%  '
$ws.Range("E37").Value = '  +1.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.82'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.16%  '
$ws.Range("E39").Value = '  +6.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0316'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.40'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.84'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.837.53'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("B46").Value = 'BitcoinSV'
$ws.Range("C46").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '82.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.94%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.197'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '75.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.88%  '
$ws.Range("E49").Value = '  +2.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '98.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '55.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.76%  '
